$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.783.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "'2.047.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'227.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").Value = "'60.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("E10").Value = "  +2.66%  "
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "'2.350.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "'21.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.15%  "
$ws.Range("D15").Value = "'5.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.64%  "
$ws.Range("D16").Value = "'0.765"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "'2.041.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "'37.717.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "'69.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "'5.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").Value = "'222.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").Value = "'2.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.68%  "
$ws.Range("D26").Value = "'169.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("D27").Value = "'9.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").Value = "'18.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  +8.13%  "
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").Value = "'0.0603"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("E37").Value = "  +4.62%  "
$ws.Range("E38").Value = "  +7.87%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").Value = "'18.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.97%  "
$ws.Range("D41").Value = "'1.531.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").Value = "'97.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "'4.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "'0.0889"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D50").Value = "'7.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").Value = "'2.238.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.01%  "
